$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Batam"
$ws.Range("B3").Value = "Ringroad"
$ws.Range("B4").Value = "Sepinggan"

$ws.Range("A5").Value = "SOMBER"
$ws.Range("B5").Value = "Batam"

$ws.Range("A6").Value = "SEPINGGAN 1"
$ws.Range("B6").Value = "Sepinggan"

$ws.Range("B7").Select()
